$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Dr. Edward Roualdes"
$ws.Range("C5").Select()
